# Add a "comment" header field to every tab of the naming template.
#
# For each worksheet in the workbook, write "comment" into the cell that
# immediately follows the existing header row (column C, since every sheet
# currently only uses columns A and B), then leave the selection sitting on
# the row below the new header (C2) - mirroring how a person would type the
# header and press Enter.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("C1").Value = "comment"
}

# Re-select C2 on every sheet (matches the header-entry workflow), then make
# sure the sheet that was active before the edit (the last tab) stays active.
foreach ($ws in $wb.Worksheets) {
    $ws.Activate()
    $ws.Range("C2").Select()
}

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
